$wb = $excel.ActiveWorkbook

# --- Sheet 1: summary_counts ---
$ws1 = $wb.Worksheets.Item("summary_counts")
$ws1.Range("B2").Value = 644
$ws1.Range("B3").Value = 251
$ws1.Range("B5").Value = 675
$ws1.Range("B6").Value = 583
$ws1.Range("B7").Value = 827
$ws1.Range("B8").Value = 568
$ws1.Range("B9").Value = 685
$ws1.Range("B10").Value = 1314
$ws1.Range("B11").Value = 680
$ws1.Range("B12").Value = 1702
$ws1.Range("B13").Value = 3705
$ws1.Range("B14").Value = 1702

# --- Sheet 2: response_action_counts ---
$ws2 = $wb.Worksheets.Item("response_action_counts")
$ws2.Range("E2").Value = 0.005548363232846311
$ws2.Range("E3").Value = 0.000369890882189754
$ws2.Range("E4").Value = 0.02126872572591086
$ws2.Range("E5").Value = 0.001294618087664139
$ws2.Range("D6").Value = 99
$ws2.Range("E6").Value = 0.01830959866839282
$ws2.Range("E7").Value = 0.002404290734233401
$ws2.Range("E8").Value = 0.005548363232846311
$ws2.Range("E9").Value = 0.000554836323284631
$ws2.Range("E10").Value = 0.006103199556130942
$ws2.Range("E11").Value = 0.001109672646569262
$ws2.Range("E12").Value = 0.004438690586277048
$ws2.Range("E13").Value = 0.001294618087664139
$ws2.Range("E14").Value = 0.01350101719992602
$ws2.Range("D15").Value = 123
$ws2.Range("E15").Value = 0.02274828925466987
$ws2.Range("E16").Value = 0.0595524320325504
$ws2.Range("E17").Value = 0.0009247272054743851
$ws2.Range("D18").Value = 24
$ws2.Range("E18").Value = 0.004438690586277048
$ws2.Range("D19").Value = 52
$ws2.Range("E19").Value = 0.009617162936933604
$ws2.Range("E20").Value = 0.003144072498612909
$ws2.Range("E21").Value = 0.000554836323284631
$ws2.Range("D22").Value = 72
$ws2.Range("E22").Value = 0.01331607175883114
$ws2.Range("E23").Value = 0.03014610689846495
$ws2.Range("E24").Value = 0.005733308673941187
$ws2.Range("E25").Value = 0.009617162936933604
$ws2.Range("E26").Value = 0.002959127057518032
$ws2.Range("D27").Value = 199
$ws2.Range("E27").Value = 0.03680414277788052
$ws2.Range("E28").Value = 0.01350101719992602
$ws2.Range("D29").Value = 17
$ws2.Range("E29").Value = 0.003144072498612909
$ws2.Range("D30").Value = 3
$ws2.Range("E30").Value = 0.000554836323284631
$ws2.Range("D31").Value = 391
$ws2.Range("E31").Value = 0.07231366746809691
$ws2.Range("E32").Value = 0.02293323469576475
$ws2.Range("D33").Value = 325
$ws2.Range("E33").Value = 0.06010726835583503
$ws2.Range("E34").Value = 0.003883854262992417
$ws2.Range("D35").Value = 35
$ws2.Range("E35").Value = 0.006473090438320696
$ws2.Range("E36").Value = 0.003513963380802663
$ws2.Range("E37").Value = 0.001294618087664139
$ws2.Range("D38").Value = 79
$ws2.Range("E38").Value = 0.01461068984649528
$ws2.Range("D39").Value = 30
$ws2.Range("E39").Value = 0.005548363232846311
$ws2.Range("E40").Value = 0.009062326613648974
$ws2.Range("D41").Value = 154
$ws2.Range("E41").Value = 0.02848159792861106
$ws2.Range("D42").Value = 8
$ws2.Range("E42").Value = 0.001479563528759016
$ws2.Range("D43").Value = 4
$ws2.Range("E43").Value = 0.0007397817643795081
$ws2.Range("E44").Value = 0.002219345293138524
$ws2.Range("E45").Value = 0.0160902533752543
$ws2.Range("D46").Value = 17
$ws2.Range("E46").Value = 0.003144072498612909
$ws2.Range("E47").Value = 0.0160902533752543
$ws2.Range("E48").Value = 0.01276123543554651
$ws2.Range("D49").Value = 53
$ws2.Range("E49").Value = 0.009802108378028482
$ws2.Range("D50").Value = 24
$ws2.Range("E50").Value = 0.004438690586277048
$ws2.Range("E51").Value = 0.007952653967079712
$ws2.Range("D52").Value = 378
$ws2.Range("E52").Value = 0.06990937673386351
$ws2.Range("D53").Value = 295
$ws2.Range("E53").Value = 0.05455890512298872
$ws2.Range("D54").Value = 277
$ws2.Range("E54").Value = 0.05122988718328093
$ws2.Range("D55").Value = 26
$ws2.Range("E55").Value = 0.004808581468466802
$ws2.Range("D56").Value = 29
$ws2.Range("E56").Value = 0.005363417791751433
$ws2.Range("D57").Value = 161
$ws2.Range("E57").Value = 0.0297762160162752
$ws2.Range("E58").Value = 0.000554836323284631
$ws2.Range("D59").Value = 104
$ws2.Range("E59").Value = 0.01923432587386721
$ws2.Range("D60").Value = 27
$ws2.Range("E60").Value = 0.004993526909561679
$ws2.Range("D61").Value = 52
$ws2.Range("E61").Value = 0.009617162936933604
$ws2.Range("D62").Value = 144
$ws2.Range("E62").Value = 0.02663214351766229
$ws2.Range("E63").Value = 0.01257628999445164
$ws2.Range("E64").Value = 0.002219345293138524
$ws2.Range("D65").Value = 60
$ws2.Range("E65").Value = 0.01109672646569262
$ws2.Range("D66").Value = 176
$ws2.Range("E66").Value = 0.03255039763269835
$ws2.Range("D67").Value = 50
$ws2.Range("E67").Value = 0.009247272054743851
$ws2.Range("E68").Value = 0.00369890882189754
$ws2.Range("D69").Value = 149
$ws2.Range("E69").Value = 0.02755687072313668
$ws2.Range("E70").Value = 0.0009247272054743851
$ws2.Range("D71").Value = 156
$ws2.Range("E71").Value = 0.02885148881080081
$ws2.Range("E72").Value = 0.01664508969853893
$ws2.Range("D73").Value = 159
$ws2.Range("E73").Value = 0.02940632513408544
$ws2.Range("D74").Value = 254
$ws2.Range("E74").Value = 0.04697614203809876
$ws2.Range("E75").Value = 0.001109672646569262
$ws2.Range("D76").Value = 105
$ws2.Range("E76").Value = 0.01941927131496209
$ws2.Range("E77").Value = 0.006103199556130942
$ws2.Range("E78").Value = 0.000554836323284631
$ws2.Range("D79").Value = 20
$ws2.Range("E79").Value = 0.00369890882189754
$ws2.Range("D80").Value = 23
$ws2.Range("E80").Value = 0.004253745145182172
$ws2.Range("D81").Value = 1702
